# Populate the "estadisticas_ofertas" summary sheet with the offer statistics
# table: a bold/centered/bordered header row followed by an offer detail row
# and a totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:E1) -----------------------------------------------
# Build the combined header style (bold font + thin box border + centered
# horizontal / top vertical alignment) on a single cell first so only one
# extra font/border/cellXf triple gets created, then fan that exact style
# out to the rest of the header via copy/paste-format (no extra styles).
$headerSeed = $ws.Range("A1")
$headerSeed.Font.Bold = $true
$headerSeed.Borders.LineStyle = 1
$headerSeed.HorizontalAlignment = -4108
$headerSeed.VerticalAlignment = -4160

$headerSeed.Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)

$ws.Range("A1").Value = "TIPO"
$ws.Range("B1").Value = "IDENTIFICADOR"
$ws.Range("C1").Value = "TOTAL ASIGNADO (kWh)"
$ws.Range("D1").Value = "PRECIO PROMEDIO"
$ws.Range("E1").Value = "COSTO TOTAL"

# --- Offer detail row (row 2) ------------------------------------------
$ws.Range("A2").Value = "OFERTA"
$ws.Range("B2").Value = "OP1_Wide -AES"
$ws.Range("C2").Value = 415000000.0000029
$ws.Range("D2").Value = 135.2911123114241
$ws.Range("E2").Value = 56145811609.24137

# --- Totals row (row 3) -------------------------------------------------
$ws.Range("A3").Value = "TOTAL"
$ws.Range("B3").Value = "TODAS LAS OFERTAS"
$ws.Range("C3").Value = 415000000.0000029
$ws.Range("D3").Value = 135.2911123114241
$ws.Range("E3").Value = 56145811609.24137
